$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.769.68"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.895.15"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7617"
$ws.Range("E5").Value = "  +3.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.37"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "1.894.44"
$ws.Range("E8").Value = "  +0.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3047"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.37"
$ws.Range("E10").Value = "  -3.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06816"
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07962"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "1.894.95"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7355"
$ws.Range("E14").Value = "  -4.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.142"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.71"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "29.782.42"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.82"
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.905"
$ws.Range("E19").Value = "  +2.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.90"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007674"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.905"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.18"
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.186"
$ws.Range("E26").Value = "  -1.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.61"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1289"
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.017"
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.405"
$ws.Range("E30").Value = "  +3.93%  "
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.248"
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.063"
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05213"
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.246"
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7238"
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.716"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01911"
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.773"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.129"
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4394"
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.53"
$ws.Range("E42").Value = "  -3.52%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8293"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.876"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.574"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.73"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.727"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").Value = "2.055.26"
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.89"
$ws.Range("E50").Value = "  -2.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05931"
$ws.Range("E51").Value = "  -0.19%  "
